$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("regimen-2-digitos") is re-curated from a dimension to a measure,
# matching the pattern already used by the other iaest-measure columns.
$ws.Range("B2").Value = "iaest-measure:regimen-2-digitos"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("B5").Clear()

# Column L ("direccion-provincial-nombre") is likewise re-curated from the
# sdmx-dimension:refArea dimension to its own iaest-measure.
$ws.Range("L2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
